# Add a "Save" column (H) to the s_vals sheet, mirroring the style of the
# existing header row and filling in the per-row save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same look as the other header cells (bold, centered,
# bordered) - copy the format from G1 so it reuses the identical style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Per-row "Save" values for rows 2-9 (plain numeric cells, matching the
# unstyled look of the other data columns).
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 1
    6 = 1
    7 = 0
    8 = 0
    9 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}

$excel.CutCopyMode = 0
